$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 92: Whinier than the Sword / Enchanted Koppranickel Ink (itemID 19901)
$ws.Range("H92").Value = 1220.5264
$ws.Range("I92").Value = 1341.7646
$ws.Range("J92").Value = 190
$ws.Range("K92").Value = 1341.7646
$ws.Range("L92").Value = 190
$ws.Range("M92").Value = -93.76459999999997
$ws.Range("N92").Value = -2686

# ALC row 112: Making Ends Meet / Superior Spiritbond Potion (itemID 27960)
$ws.Range("H112").Value = 1796.7894
$ws.Range("I112").Value = 1412.5
$ws.Range("J112").Value = 1899.2667
$ws.Range("K112").Value = 4237.5
$ws.Range("L112").Value = 5697.800099999999
$ws.Range("M112").Value = -3129.5
$ws.Range("N112").Value = -7913.800099999999

# ALC row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone (itemID 44013)
$ws.Range("H137").Value = 653.96295
$ws.Range("I137").Value = 582.10254
$ws.Range("J137").Value = 840.8
$ws.Range("K137").Value = 1746.30762
$ws.Range("L137").Value = 2522.4
$ws.Range("M137").Value = 803.69238
$ws.Range("N137").Value = -7622.4

# ALC row 138: All-night Crafting / Cunning Craftsman's Tisane (itemID 44169)
$ws.Range("H138").Value = 1316.6262
$ws.Range("I138").Value = 684.0784
$ws.Range("J138").Value = 1988.7084
$ws.Range("K138").Value = 2052.2352
$ws.Range("L138").Value = 5966.1252
$ws.Range("M138").Value = 3087.7648
$ws.Range("N138").Value = -16246.1252

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32: Ingot We Trust / Steel Ingot (itemID 44147)
$ws.Range("H32").Value = 3332.6235
$ws.Range("I32").Value = 2686.671
$ws.Range("J32").Value = 6977.643
$ws.Range("K32").Value = 2686.671
$ws.Range("L32").Value = 6977.643
$ws.Range("M32").Value = -2399.671
$ws.Range("N32").Value = -7551.643

# ARM row 61: Dealing with the Tough Stuff / Cobalt Ingot (itemID 43999)
$ws.Range("H61").Value = 885
$ws.Range("I61").Value = 811.08105
$ws.Range("J61").Value = 1796.6666
$ws.Range("K61").Value = 811.08105
$ws.Range("L61").Value = 1796.6666
$ws.Range("M61").Value = -599.08105
$ws.Range("N61").Value = -2220.6666

# ARM row 74: As the Bolt Flies / Titanium Nugget (itemID 44000)
$ws.Range("H74").Value = 1335.5714
$ws.Range("I74").Value = 1368
$ws.Range("J74").Value = 1216.6666
$ws.Range("K74").Value = 1368
$ws.Range("L74").Value = 1216.6666
$ws.Range("M74").Value = -494
$ws.Range("N74").Value = -2964.6666

# ARM row 77: Heavy Metal Banned (L) / Titanium Nugget (itemID 44000)
$ws.Range("H77").Value = 1335.5714
$ws.Range("I77").Value = 1368
$ws.Range("J77").Value = 1216.6666
$ws.Range("K77").Value = 6840
$ws.Range("L77").Value = 6083.333000000001
$ws.Range("M77").Value = -2472
$ws.Range("N77").Value = -14819.333

# ARM row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot (itemID 43997)
$ws.Range("H132").Value = 2288.7693
$ws.Range("I132").Value = 2104.5417
$ws.Range("K132").Value = 6313.625100000001
$ws.Range("M132").Value = -3783.625100000001

# ARM row 136: Metal with Mettle / Cobalt Tungsten Ingot (itemID 43999)
$ws.Range("H136").Value = 885
$ws.Range("I136").Value = 811.08105
$ws.Range("J136").Value = 1796.6666
$ws.Range("K136").Value = 2433.24315
$ws.Range("L136").Value = 5389.9998
$ws.Range("M136").Value = 116.7568499999998
$ws.Range("N136").Value = -10489.9998

$ws = $wb.Worksheets.Item("BSM")
# BSM row 86: Through Thick and Thin / Adamantite Nugget (itemID 12526)
$ws.Range("H86").Value = 2254.9546
$ws.Range("I86").Value = 1981.381
$ws.Range("K86").Value = 1981.381
$ws.Range("M86").Value = -858.3810000000001

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget (itemID 12526)
$ws.Range("H89").Value = 2254.9546
$ws.Range("I89").Value = 1981.381
$ws.Range("K89").Value = 9906.905000000001
$ws.Range("M89").Value = -4290.905000000001

# BSM row 134: Ruthenium Supremium / Ruthenium Ingot (itemID 43998)
$ws.Range("H134").Value = 19563.457
$ws.Range("I134").Value = 1773.5319
$ws.Range("J134").Value = 103176.1
$ws.Range("K134").Value = 5320.5957
$ws.Range("L134").Value = 309528.3
$ws.Range("M134").Value = -2785.5957
$ws.Range("N134").Value = -314598.3

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31: Wall Not Found / Walnut Lumber (itemID 44023)
$ws.Range("H31").Value = 1668871.2
$ws.Range("I31").Value = 2144947.2
$ws.Range("J31").Value = 2605.3333
$ws.Range("K31").Value = 2144947.2
$ws.Range("L31").Value = 2605.3333
$ws.Range("M31").Value = -2144652.2
$ws.Range("N31").Value = -3195.3333

# CRP row 34: Armoires of the Rich and Famous / Walnut Lumber (itemID 44023)
$ws.Range("H34").Value = 1668871.2
$ws.Range("I34").Value = 2144947.2
$ws.Range("J34").Value = 2605.3333
$ws.Range("K34").Value = 2144947.2
$ws.Range("L34").Value = 2605.3333
$ws.Range("M34").Value = -2144745.2
$ws.Range("N34").Value = -3009.3333

# CRP row 58: You Do the Heavy Lifting / Mahogany Lumber (itemID 44021)
$ws.Range("H58").Value = 5686.6
$ws.Range("I58").Value = 1911.6666
$ws.Range("J58").Value = 11349
$ws.Range("K58").Value = 1911.6666
$ws.Range("L58").Value = 11349
$ws.Range("M58").Value = -1708.6666
$ws.Range("N58").Value = -11755

# CRP row 132: Hull Lotta Damage / Ginseng Lumber (itemID 44019)
$ws.Range("H132").Value = 1621.9436
$ws.Range("I132").Value = 995.9487
$ws.Range("J132").Value = 2384.875
$ws.Range("K132").Value = 2987.8461
$ws.Range("L132").Value = 7154.625
$ws.Range("M132").Value = -457.8461000000002
$ws.Range("N132").Value = -12214.625

# CRP row 134: Wood You Be Quiet / Ceiba Lumber (itemID 44020)
$ws.Range("H134").Value = 1526.8723
$ws.Range("I134").Value = 1517.1562
$ws.Range("J134").Value = 1547.6
$ws.Range("K134").Value = 4551.4686
$ws.Range("L134").Value = 4642.799999999999
$ws.Range("M134").Value = -2016.4686
$ws.Range("N134").Value = -9712.799999999999

# CRP row 136: Turali Quality / Dark Mahogany Lumber (itemID 44021)
$ws.Range("H136").Value = 5686.6
$ws.Range("I136").Value = 1911.6666
$ws.Range("J136").Value = 11349
$ws.Range("K136").Value = 5734.9998
$ws.Range("L136").Value = 34047
$ws.Range("M136").Value = -3184.9998
$ws.Range("N136").Value = -39147

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5: What a Sap / Maple Syrup (itemID 43974)
$ws.Range("H5").Value = 446.21875
$ws.Range("I5").Value = 303.5
$ws.Range("J5").Value = 874.375
$ws.Range("K5").Value = 910.5
$ws.Range("L5").Value = 2623.125
$ws.Range("M5").Value = -798.5
$ws.Range("N5").Value = -2847.125

# CUL row 113: Can't Eat Just One / Night Vinegar (itemID 27843)
$ws.Range("H113").Value = 560.9545000000001
$ws.Range("I113").Value = 1407.5
$ws.Range("J113").Value = 476.3
$ws.Range("K113").Value = 4222.5
$ws.Range("L113").Value = 1428.9
$ws.Range("M113").Value = -2052.5
$ws.Range("N113").Value = -5768.9

# CUL row 122: Salt of the North / Northern Sea Salt (itemID 36078)
$ws.Range("H122").Value = 620.1
$ws.Range("J122").Value = 635.8889
$ws.Range("L122").Value = 5723.0001
$ws.Range("N122").Value = -10623.0001

# CUL row 131: The Mountain Steeped / Tsai tou Vounou (itemID 36060)
$ws.Range("H131").Value = 32290.39
$ws.Range("I131").Value = 202150
$ws.Range("J131").Value = 17895.508
$ws.Range("K131").Value = 606450
$ws.Range("L131").Value = 53686.524
$ws.Range("M131").Value = -601410
$ws.Range("N131").Value = -63766.524

# CUL row 135: Not-so-secret Ingredient / Royal Maple Syrup (itemID 43974)
$ws.Range("H135").Value = 446.21875
$ws.Range("I135").Value = 303.5
$ws.Range("J135").Value = 874.375
$ws.Range("K135").Value = 2731.5
$ws.Range("L135").Value = 7869.375
$ws.Range("M135").Value = -196.5
$ws.Range("N135").Value = -12939.375

# CUL row 140: Sweet, Sweet Bean Juice / Mesquite Juice (itemID 44097)
$ws.Range("H140").Value = 105794.555
$ws.Range("I140").Value = 132136.22
$ws.Range("K140").Value = 396408.66
$ws.Range("M140").Value = -391228.66

$ws = $wb.Worksheets.Item("GSM")
# GSM row 53: North Ore South / Electrum Gorget (itemID 4361)
$ws.Range("H53").Value = 7167.6
$ws.Range("I53").Value = 7959.3335
$ws.Range("J53").Value = 5980
$ws.Range("K53").Value = 7959.3335
$ws.Range("L53").Value = 5980
$ws.Range("M53").Value = -7328.3335
$ws.Range("N53").Value = -7242

# GSM row 58: The Big Red / Red Coral Necklace (itemID 4363)
$ws.Range("H58").Value = 8000
$ws.Range("I58").Value = 6000
$ws.Range("K58").Value = 6000
$ws.Range("M58").Value = -5723

# GSM row 107: Whetstones for the Workers / Hard Mudstone Whetstone (itemID 27802)
$ws.Range("H107").Value = 172
$ws.Range("I107").Value = 172
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 172
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 1748

# GSM row 132: On Board for Lar / Lar Ingot (itemID 44008)
$ws.Range("H132").Value = 2024.3541
$ws.Range("I132").Value = 1889.8572
$ws.Range("J132").Value = 2212.65
$ws.Range("K132").Value = 5669.571599999999
$ws.Range("L132").Value = 6637.950000000001
$ws.Range("M132").Value = -3139.571599999999
$ws.Range("N132").Value = -11697.95

$ws = $wb.Worksheets.Item("LTW")
# LTW row 93: Hide to Go Seek / Gagana Leather (itemID 19993)
$ws.Range("H93").Value = 2800
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

# LTW row 132: Tenets of Tanning / Silver Lobo Leather (itemID 44058)
$ws.Range("H132").Value = 2076.0386
$ws.Range("I132").Value = 1888.5897
$ws.Range("J132").Value = 2638.3845
$ws.Range("K132").Value = 5665.7691
$ws.Range("L132").Value = 7915.1535
$ws.Range("M132").Value = -3135.7691
$ws.Range("N132").Value = -12975.1535

$ws = $wb.Worksheets.Item("WVR")
# WVR row 96: Skills on Display / Ruby Cotton Cloth (itemID 19977)
$ws.Range("H96").Value = 2600.75
$ws.Range("I96").Value = 2600.75
$ws.Range("K96").Value = 2600.75
$ws.Range("M96").Value = -1227.75

# WVR row 132: Comfy Cabins / Snow Cotton Cloth (itemID 44029)
$ws.Range("H132").Value = 1854.5714
$ws.Range("I132").Value = 1596
$ws.Range("J132").Value = 2501
$ws.Range("K132").Value = 4788
$ws.Range("L132").Value = 7503
$ws.Range("M132").Value = -2258
$ws.Range("N132").Value = -12563

# WVR row 136: Weaving the Envelope / Sarcenet Cloth (itemID 44031)
$ws.Range("H136").Value = 1774.4117
$ws.Range("I136").Value = 2740
$ws.Range("J136").Value = 688.125
$ws.Range("K136").Value = 8220
$ws.Range("L136").Value = 2064.375
$ws.Range("M136").Value = -5670
$ws.Range("N136").Value = -7164.375
